$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

foreach ($r in 8..15) {
    $ws.Range("E$r").Value = "z"
}
$ws.Range("E17").Value = "z"

$ws.Application.ActiveWindow.Zoom = 115
$ws.Range("M14").Select()
